$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 (I0) and J1 (IF)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header style/formatting from H1 (bold, bordered, centered) onto I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data for columns I (I0) and J (IF) per row: @(row, I-value, J-value)
$data = @(
    @(2,9,9),
    @(3,9,9),
    @(4,2,2),
    @(5,6,6),
    @(6,9,9),
    @(7,8,9),
    @(8,7,7),
    @(9,7,7),
    @(10,6,7),
    @(11,10,10),
    @(12,8,8),
    @(13,7,7),
    @(14,7,7),
    @(15,8,8),
    @(16,8,8),
    @(17,7,7),
    @(18,8,8),
    @(19,7,7),
    @(20,9,9),
    @(21,8,8),
    @(22,7,7),
    @(23,8,8),
    @(24,9,9),
    @(25,7,7),
    @(26,5,5),
    @(27,8,8),
    @(28,9,9),
    @(29,6,6),
    @(30,8,8),
    @(31,7,7),
    @(32,8,8),
    @(33,7,7),
    @(34,6,6),
    @(35,5,6),
    @(36,6,6),
    @(37,9,9),
    @(38,8,8),
    @(39,7,7),
    @(40,8,8),
    @(41,7,8),
    @(42,7,8),
    @(43,7,8),
    @(44,9,9),
    @(45,9,9),
    @(46,8,8),
    @(47,6,6),
    @(48,9,9),
    @(49,8,8),
    @(50,7,7),
    @(51,5,5),
    @(52,10,10),
    @(53,7,8),
    @(54,6,6),
    @(55,7,7),
    @(56,8,8),
    @(57,6,6),
    @(58,7,7),
    @(59,7,7),
    @(60,8,8),
    @(61,7,7),
    @(62,9,9),
    @(63,7,7),
    @(64,9,9),
    @(65,9,9),
    @(66,8,8),
    @(67,7,7),
    @(68,8,8),
    @(69,7,7),
    @(70,8,8),
    @(71,6,7),
    @(72,8,8),
    @(73,8,8),
    @(74,9,9),
    @(75,7,7),
    @(76,8,8),
    @(77,6,7),
    @(78,6,6),
    @(79,7,8),
    @(80,7,7),
    @(81,5,5),
    @(82,5,5),
    @(83,8,8),
    @(84,4,5),
    @(85,4,4)
)

foreach ($entry in $data) {
    $row = $entry[0]
    $ws.Cells.Item($row, 9).Value = $entry[1]
    $ws.Cells.Item($row, 10).Value = $entry[2]
}
